# Add a new "dipwell_measurements" pointer row (comparing modelled to
# measured dipwell data) below the existing "initial_dipwell_measurements"
# row on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A16").Value = "dipwell_measurements"
$ws.Range("B16").Value = "data/dipwell_data_from_first_rainfall_record.csv"

# Match the author's resulting selection/active cell.
$ws.Range("A16").Select()
